$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# Sheet1: insert a new "Status" column before the current column C
# ---------------------------------------------------------------------------
$ws1.Columns.Item(3).Insert()
$ws1.Columns.Item(3).ColumnWidth = 15.65

# Header
$ws1.Range("C1").Value = "Status"

# Rename Temple -> Chapel (row 16, now in column D after the insert)
$ws1.Range("D16").Value = "Chapel"

# "in Engine" status markers
$ws1.Range("C2").Value = "in Engine"
$ws1.Range("C2").HorizontalAlignment = -4131

$ws1.Range("C10").Value = "in Engine"
$ws1.Range("C10").HorizontalAlignment = -4131
$ws1.Range("C10").VerticalAlignment = -4108

$ws1.Range("C9").HorizontalAlignment = -4131
$ws1.Range("C9").VerticalAlignment = -4108
$ws1.Range("C11").HorizontalAlignment = -4131
$ws1.Range("C11").VerticalAlignment = -4108
$ws1.Range("C12").HorizontalAlignment = -4131
$ws1.Range("C12").VerticalAlignment = -4108
$ws1.Range("C13").HorizontalAlignment = -4131
$ws1.Range("C13").VerticalAlignment = -4108

$ws1.Range("C16").Value = "in Engine"
$ws1.Range("C16").HorizontalAlignment = -4131

$ws1.Range("C17").Value = "in Engine"
$ws1.Range("C17").HorizontalAlignment = -4131

$ws1.Range("C44").Value = "in Engine"
$ws1.Range("C44").HorizontalAlignment = -4131

$ws1.Range("C46").Value = "in Engine"
$ws1.Range("C46").HorizontalAlignment = -4131

# New building rows
$ws1.Range("D50").Value = "Courthouse"
$ws1.Range("C50").Value = "in Engine"
$ws1.Range("C50").HorizontalAlignment = -4131

$ws1.Range("D51").Value = "Town Hall"
$ws1.Range("D52").Value = "Palace"

# Data validation dropdown on the whole Status column
$ws1.Range("C2:C52").Validation.Add(3, 1, 1, '"Not started, In progress, 3D done, in Engine"')

# Selection on Sheet1
$ws1.Activate() | Out-Null
$ws1.Range("D19").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet2 cosmetic tweaks
# ---------------------------------------------------------------------------
$ws2.Activate() | Out-Null
$ws2.Columns.Item(1).ColumnWidth = 13.65
$ws2.Range("E9").Select() | Out-Null

# restore Sheet1 as the active sheet/tab
$ws1.Activate() | Out-Null
$ws1.Range("D19").Select() | Out-Null
